{"js": "async (context) => {\n  // The underlying change in this revision only touches the auto-generated\n  // XML namespace-prefix numbering that Word stamps on the root element of\n  // several package parts (document.xml, headers/footers, styles.xml,\n  // numbering.xml, theme1.xml, ...) when it re-saves the file. No visible\n  // text, formatting, structure, or content actually changed between the\n  // two revisions (the \"ns8\"/\"ns19\" aliases that became \"ns9\"/\"ns17\" are\n  // declared but never referenced anywhere in the markup). The commit\n  // message describes an application/back-end fix unrelated to this\n  // document's content.\n  //\n  // Touch the body (load/sync) so the context is exercised, but make no\n  // content changes - this mirrors the no-semantic-change nature of the\n  // diff.\n  const body = context.document.body;\n  body.load(\"text\");\n  await context.sync();\n};\n", "ps1": "# The captured revision only differs in the auto-generated XML namespace-\n# prefix numbering Word stamps on the root element of several package parts\n# (document.xml, header/footer/footnote/endnote parts, styles.xml,\n# numbering.xml, theme1.xml, ...) when the file is re-saved - e.g. \"ns8\"/\n# \"ns19\" become \"ns9\"/\"ns17\". Those prefixes are declared but never actually\n# referenced anywhere in the markup, so no visible text, formatting,\n# structure, or content changed between the two revisions. The commit\n# message itself describes an unrelated application/back-end bug fix (a\n# company-registration screen's contact field and a boolean-to-integer type\n# change), not an edit to this document's content.\n#\n# Touch the document (read-only) so the automation path is exercised, but\n# make no content changes - this mirrors the no-semantic-change nature of\n# the diff.\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
